$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("13_03_2019")
$ws.Activate()

# 1) Flip the "COMPLETO" answer for the "Adicionar en verificacion..." row from NO to SI
$ws.Range("C2").Value = "SI"

# 2) Add the new observaciones block below the existing table (rows 15-18)
$ws.Range("A15").Value = "REQUERIMIENTOS"
$ws.Range("B15").Value = "OBSERVACIONES"

$ws.Range("A16").Value = "Adicionar en verificación preguntas de si o no."
$ws.Range("B16").Value = "El apartado de INFORMACIÓN FINANCIERA, ACCIONISTAS y DECLARACIÓN ORIGEN DE FONDOS en JURIDICO no se muestra cuando el estado del formulario es VERIFICACION, teniendo en cuenta que este contiene campos que necesitan verificarse."

$ws.Range("A17").Value = "Adicionar en verificación preguntas de si o no."
$ws.Range("B17").Value = "El campos CIIU no existe en el módulo de DATOS EMPRESA DONDE TRABAJA en NATURAL"

$ws.Range("A18").Value = "Adicionar en verificación preguntas de si o no."
$ws.Range("B18").Value = "El campo declaración origen de fondos no existe dentro del apartado de INFORMACIÓN FINANCIERA"

# 3) Turn A15:B18 into a proper table ("Tabla3"), matching Tabla136/Tabla2 style
$tbl = $ws.ListObjects.Add(1, $ws.Range("A15:B18"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Tabla3"
$tbl.TableStyle = "TableStyleMedium2"

# 4) Update view: zoom out to 70% and move selection to C18
$excel.ActiveWindow.Zoom = 70
$ws.Range("C18").Select()
